# Daily attendance processing - 2025-10-31 08:27:04
# Rotates the "Recorded By" (column G) comma-separated list left by one
# position whenever the first entry is not "System"/"system". This moves
# a leading non-System identity (e.g. an email address) to the end of
# the list, leaving lists that already start with System untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1 -and $parts[0].Trim().ToLower() -ne "system") {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
